# Update countries & provincias Spain
# (this workbook covers the "Pais" country data sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 12:35"

# --- Plain data refresh (no country re-sort) ---------------------------
# Alemania (row 11)
$ws.Range("B11").Value = 173274
$ws.Range("C11").Value = 103
$ws.Range("G11").Value = 17
$ws.Range("H11").Value = 7755

# Suiza (row 25)
$ws.Range("B25").Value = 30413
$ws.Range("C25").Value = 33
$ws.Range("E25").Value = 1746

# Australia (row 53)
$ws.Range("D53").Value = 6271
$ws.Range("E53").Value = 611
$ws.Range("F53").Value = 17

# Marruecos (row 56)
$ws.Range("B56").Value = 6466
$ws.Range("C56").Value = 48
$ws.Range("D56").Value = 3099
$ws.Range("E56").Value = 3179

# Sri Lanka (row 104)
$ws.Range("B104").Value = 891
$ws.Range("C104").Value = 2
$ws.Range("E104").Value = 500

# Albania (row 105)
$ws.Range("B105").Value = 880
$ws.Range("C105").Value = 4
$ws.Range("D105").Value = 688
$ws.Range("E105").Value = 161

# --- Countries re-sorted: Afganistan now ranks above Moldavia/Ghana ----
# Row 61 becomes Afganistan with refreshed totals
$ws.Range("A61").Value = "Afganistan"
$ws.Range("B61").Value = 5226
$ws.Range("C61").Value = 263
$ws.Range("D61").Value = 648
$ws.Range("E61").Value = 4446
$ws.Range("F61").Value = 7
$ws.Range("G61").Value = 5
$ws.Range("H61").Value = 132

# Row 62 becomes Moldavia (carries the old Afganistan row's figures)
$ws.Range("A62").Value = "Moldavia"
$ws.Range("B62").Value = 5154
$ws.Range("C62").Value = 0
$ws.Range("D62").Value = 2069
$ws.Range("E62").Value = 2903
$ws.Range("F62").Value = 251
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 182

# Row 63 becomes Ghana (carries the old Moldavia row's figures)
$ws.Range("A63").Value = "Ghana"
$ws.Range("B63").Value = 5127
$ws.Range("C63").Value = 0
$ws.Range("D63").Value = 494
$ws.Range("E63").Value = 4611
$ws.Range("F63").Value = 5
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 22

# --- Curazao / Dominica swap (row 198 / 199) ---------------------------
$ws.Range("A198").Value = "Dominica"
$ws.Range("D198").Value = 15
$ws.Range("H198").Value = 0

$ws.Range("A199").Value = "Curazao"
$ws.Range("D199").Value = 14
$ws.Range("H199").Value = 1

# --- San Bartolome / Bonaire swap (row 215 / 216, figures tied) -------
$ws.Range("A215").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A216").Value = "San Bartolome"
